$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Rows 757-766 of the Translation sheet were blank placeholder rows; they
# now get populated with 10 new text entries (TEXT ID / TYPOGRAPHY NAME /
# ALIGNMENT / DIRECTION / GB / IT / ES / FR), alternating between a
# "Center" alignment + literal "<value>" placeholder, and a "Left"
# alignment + literal "0" placeholder.
$data = @(
    @("SingleUseId4166", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId4167", "Default", "Left",   "LTR", "0"),
    @("SingleUseId4168", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId4169", "Default", "Left",   "LTR", "0"),
    @("SingleUseId4170", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId4171", "Default", "Left",   "LTR", "0"),
    @("SingleUseId4172", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId4173", "Default", "Left",   "LTR", "0"),
    @("SingleUseId4174", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId4175", "Default", "Left",   "LTR", "0")
)

$startRow = 757
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]

    $valueText = $vals[4]
    $langRange = $ws.Range($ws.Cells.Item($row, 6), $ws.Cells.Item($row, 9))
    if ($valueText -eq "0") {
        # "0" has to stay text (not become the number 0). Entering it via a
        # formula that evaluates to the string "0" and then converting the
        # range to static values keeps the shared-string "text" type
        # without stamping a text number-format on the cell.
        $langRange.Formula = '="0"'
        $langRange.Copy()
        $langRange.PasteSpecial(-4163)
    } else {
        $langRange.Value = $valueText
    }
}
$excel.CutCopyMode = $false
